$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Dije Flor Rosa"
$ws.Range("E2").Value = "Aleación"
$ws.Range("F2").Value = "Rosa"
$ws.Range("G2").Value = "Mujer"
$ws.Range("H2").Value = 22.03
$ws.Range("I2").Value = "['dije', 'flor', 'rosa', 'aleación', 'joyería', 'dorado']"
$ws.Range("J2").Value = "Dije en forma de flor con detalles dorados en el centro y en los bordes."
$ws.Range("K2").Value = "Este adorable dije en forma de flor es perfecto para cualquier joyería. Con detalles dorados en el centro y en los bordes, este dije añade un toque de elegancia a cualquier collar o pulsera. El material aleación asegura un acabado duradero y resistente."
$ws.Range("L2").Value = "dije_flor_pink"
